$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 8.366847333333332
$ws.Range("H2").Value = 25.100542
$ws.Range("I2").Value = 0.3599008094639652
$ws.Range("J2").Value = 0.3608281410662089
$ws.Range("M2").Value = 38.84293433333333
$ws.Range("N2").Value = 116.528803
$ws.Range("O2").Value = 0.4759825783603507
$ws.Range("P2").Value = 0.5627959792472441
$ws.Range("Q2").Value = 324.9929015456917
$ws.Range("R2").Value = 2924.936113911225
$ws.Range("S2").Value = 0.1713065152426355
$ws.Range("T2").Value = 0.2030726269913198

$ws.Range("G3").Value = 8.366847333333332
$ws.Range("H3").Value = 25.100542
$ws.Range("I3").Value = 0.3599008094639652
$ws.Range("J3").Value = 0.3608281410662089
$ws.Range("M3").Value = 2.977970666666666
$ws.Range("N3").Value = 8.933911999999999
$ws.Range("O3").Value = 0.03649214922944397
$ws.Range("P3").Value = 0.043147870939245
$ws.Range("Q3").Value = 24.91622593114488
$ws.Range("R3").Value = 224.246033380304
$ws.Range("S3").Value = 0.0131335540467567
$ws.Range("T3").Value = 0.01556896606197247

$ws.Range("G4").Value = 8.366847333333332
$ws.Range("H4").Value = 25.100542
$ws.Range("I4").Value = 0.3599008094639652
$ws.Range("J4").Value = 0.3608281410662089
$ws.Range("M4").Value = 1.675036
$ws.Range("N4").Value = 5.025107999999999
$ws.Range("O4").Value = 0.02052594552420851
$ws.Range("P4").Value = 0.02426962694951188
$ws.Range("Q4").Value = 14.01477048983733
$ws.Range("R4").Value = 126.132934408536
$ws.Range("S4").Value = 0.007387304409175899
$ws.Range("T4").Value = 0.008757164376562735

$ws.Range("G5").Value = 8.366847333333332
$ws.Range("H5").Value = 25.100542
$ws.Range("I5").Value = 0.3599008094639652
$ws.Range("J5").Value = 0.3608281410662089
$ws.Range("M5").Value = 37.764002
$ws.Range("N5").Value = 75.52800400000001
$ws.Range("O5").Value = 0.4627613065200399
$ws.Range("P5").Value = 0.3647755394155193
$ws.Range("Q5").Value = 315.9656394296946
$ws.Range("R5").Value = 1895.793836578168
$ws.Range("S5").Value = 0.1665481688051645
$ws.Range("T5").Value = 0.1316212797937254

$ws.Range("G6").Value = 8.366847333333332
$ws.Range("H6").Value = 25.100542
$ws.Range("I6").Value = 0.3599008094639652
$ws.Range("J6").Value = 0.3608281410662089
$ws.Range("M6").Value = 0.345847
$ws.Range("N6").Value = 1.037541
$ws.Range("O6").Value = 0.004238020365956876
$ws.Range("P6").Value = 0.005010983448479815
$ws.Range("Q6").Value = 2.893649049691333
$ws.Range("R6").Value = 26.042841447222
$ws.Range("S6").Value = 0.00152526696023265
$ws.Range("T6").Value = 0.001808103842628512

$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.17924
$ws.Range("H7").Value = 0.35848
$ws.Range("I7").Value = 0.007710027268134825
$ws.Range("J7").Value = 0.005153262109217186
$ws.Range("M7").Value = 38.84293433333333
$ws.Range("N7").Value = 116.528803
$ws.Range("O7").Value = 0.4759825783603507
$ws.Range("P7").Value = 0.5627959792472441
$ws.Range("Q7").Value = 6.962207549906665
$ws.Range("R7").Value = 41.77324529944
$ws.Range("S7").Value = 0.003669838658315425
$ws.Range("T7").Value = 0.002900235195074605

$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.17924
$ws.Range("H8").Value = 0.35848
$ws.Range("I8").Value = 0.007710027268134825
$ws.Range("J8").Value = 0.005153262109217186
$ws.Range("M8").Value = 2.977970666666666
$ws.Range("N8").Value = 8.933911999999999
$ws.Range("O8").Value = 0.03649214922944397
$ws.Range("P8").Value = 0.043147870939245
$ws.Range("Q8").Value = 0.5337714622933333
$ws.Range("R8").Value = 3.20262877376
$ws.Range("S8").Value = 0.0002813554656318583
$ws.Range("T8").Value = 0.0002223522884046046

$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.17924
$ws.Range("H9").Value = 0.35848
$ws.Range("I9").Value = 0.007710027268134825
$ws.Range("J9").Value = 0.005153262109217186
$ws.Range("M9").Value = 1.675036
$ws.Range("N9").Value = 5.025107999999999
$ws.Range("O9").Value = 0.02052594552420851
$ws.Range("P9").Value = 0.02426962694951188
$ws.Range("Q9").Value = 0.30023345264
$ws.Range("R9").Value = 1.80140071584
$ws.Range("S9").Value = 0.0001582555996958976
$ws.Range("T9").Value = 0.0001250677489637558

$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.17924
$ws.Range("H10").Value = 0.35848
$ws.Range("I10").Value = 0.007710027268134825
$ws.Range("J10").Value = 0.005153262109217186
$ws.Range("M10").Value = 37.764002
$ws.Range("N10").Value = 75.52800400000001
$ws.Range("O10").Value = 0.4627613065200399
$ws.Range("P10").Value = 0.3647755394155193
$ws.Range("Q10").Value = 6.768819718480001
$ws.Range("R10").Value = 27.07527887392001
$ws.Range("S10").Value = 0.003567902291907205
$ws.Range("T10").Value = 0.001879783965639256

$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.17924
$ws.Range("H11").Value = 0.35848
$ws.Range("I11").Value = 0.007710027268134825
$ws.Range("J11").Value = 0.005153262109217186
$ws.Range("M11").Value = 0.345847
$ws.Range("N11").Value = 1.037541
$ws.Range("O11").Value = 0.004238020365956876
$ws.Range("P11").Value = 0.005010983448479815
$ws.Range("Q11").Value = 0.06198961628000001
$ws.Range("R11").Value = 0.37193769768
$ws.Range("S11").Value = 3.267525258443824 / 100000
$ws.Range("T11").Value = 2.58229111349655 / 100000

$ws.Range("G12").Value = 14.70156066666667
$ws.Range("H12").Value = 44.104682
$ws.Range("I12").Value = 0.6323891632679
$ws.Range("J12").Value = 0.6340185968245738
$ws.Range("M12").Value = 38.84293433333333
$ws.Range("N12").Value = 116.528803
$ws.Range("O12").Value = 0.4759825783603507
$ws.Range("P12").Value = 0.5627959792472441
$ws.Range("Q12").Value = 571.0517555728494
$ws.Range("R12").Value = 5139.465800155645
$ws.Range("S12").Value = 0.3010062244593998
$ws.Range("T12").Value = 0.3568231170608497

$ws.Range("G13").Value = 14.70156066666667
$ws.Range("H13").Value = 44.104682
$ws.Range("I13").Value = 0.6323891632679
$ws.Range("J13").Value = 0.6340185968245738
$ws.Range("M13").Value = 2.977970666666666
$ws.Range("N13").Value = 8.933911999999999
$ws.Range("O13").Value = 0.03649214922944397
$ws.Range("P13").Value = 0.043147870939245
$ws.Range("Q13").Value = 43.78081641955377
$ws.Range("R13").Value = 394.027347775984
$ws.Range("S13").Value = 0.02307723971705541
$ws.Range("T13").Value = 0.02735655258886792

$ws.Range("G14").Value = 14.70156066666667
$ws.Range("H14").Value = 44.104682
$ws.Range("I14").Value = 0.6323891632679
$ws.Range("J14").Value = 0.6340185968245738
$ws.Range("M14").Value = 1.675036
$ws.Range("N14").Value = 5.025107999999999
$ws.Range("O14").Value = 0.02052594552420851
$ws.Range("P14").Value = 0.02426962694951188
$ws.Range("Q14").Value = 24.62564337285066
$ws.Range("R14").Value = 221.630790355656
$ws.Range("S14").Value = 0.01298038551533672
$ws.Range("T14").Value = 0.01538739482398538

$ws.Range("G15").Value = 14.70156066666667
$ws.Range("H15").Value = 44.104682
$ws.Range("I15").Value = 0.6323891632679
$ws.Range("J15").Value = 0.6340185968245738
$ws.Range("M15").Value = 37.764002
$ws.Range("N15").Value = 75.52800400000001
$ws.Range("O15").Value = 0.4627613065200399
$ws.Range("P15").Value = 0.3647755394155193
$ws.Range("Q15").Value = 555.1897664191214
$ws.Range("R15").Value = 3331.138598514728
$ws.Range("S15").Value = 0.2926452354229682
$ws.Range("T15").Value = 0.2312744756561546

$ws.Range("G16").Value = 14.70156066666667
$ws.Range("H16").Value = 44.104682
$ws.Range("I16").Value = 0.6323891632679
$ws.Range("J16").Value = 0.6340185968245738
$ws.Range("M16").Value = 0.345847
$ws.Range("N16").Value = 1.037541
$ws.Range("O16").Value = 0.004238020365956876
$ws.Range("P16").Value = 0.005010983448479815
$ws.Range("Q16").Value = 5.084490651884667
$ws.Range("R16").Value = 45.760415866962
$ws.Range("S16").Value = 0.002680078153139788
$ws.Range("T16").Value = 0.003177056694716337
